# Add the "2022-Q4" sheet (copy of "2021-Q4" to inherit formatting/styles),
# position it right after "总计" and before "2021-Q4", then overwrite its
# values with the new fund data. Finally, insert a corresponding summary
# row into "总计".

$wb = $excel.ActiveWorkbook

# --- 1) Create the new "2022-Q4" sheet by copying "2021-Q4" -------------
$src = $wb.Worksheets.Item("2021-Q4")
$src.Copy($src)                      # new copy is placed immediately before $src
$new = $wb.Worksheets.Item(2)        # 总计, [new copy], 2021-Q4, 2020-Q4
$new.Name = "2022-Q4"

# "2021-Q4" has two data rows; 2022-Q4 only needs one, so drop row 3.
$new.Rows(3).Delete()

# --- 2) Overwrite the copied sheet's values with the 2022-Q4 fund data --
$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

$new.Range("A2").Value = 0
$new.Range("B2").Value = "'001849"
$new.Range("C2").Value = "前海开源强势共识100强等权重股票"
$new.Range("D2").Value = "'0.26"
$new.Range("E2").Value = "'71.91"
$new.Range("F2").Value = "'0.89"
$new.Range("G2").Value = "'0.0023"
$new.Range("H2").Value = 6

# --- 3) Insert the "2022-Q4" summary row into "总计" ---------------------
# Shift the existing "2021-Q4"/"2020-Q4" rows down by one row, then write
# the new "2022-Q4" row in the gap that opens up at row 2.
$tot = $wb.Worksheets.Item("总计")

$tot.Range("A4").Value = 2
$tot.Range("B4").Value = "2020-Q4"
$tot.Range("C4").Value = 2
$tot.Range("D4").Value = 0.04

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = "2021-Q4"
$tot.Range("C3").Value = 2
$tot.Range("D3").Value = 0.03

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q4"
$tot.Range("C2").Value = 1
$tot.Range("D2").Value = 0
